# RPA datasets push 2024-07-16
# Insert a new IPO record (엑셀세라퓨틱스) as row 2, shifting existing data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (push rows 2..22 down to 3..23).
$ws.Rows.Item(2).Insert()

# Row-insert copies the header row's bold/bordered formatting down onto the new
# row; strip that so the new data row matches the other plain (unstyled) rows.
$ws.Rows.Item(2).ClearFormats()

# The date-like columns (A/D/E) hold plain text such as "2024-07-03" in this
# workbook, not real dates. Mark them as Text before assigning so Excel does
# not silently reinterpret the string as a date serial, then drop back to the
# Normal style so no stray number-format style sticks to the cell.
$dateCols = "A2", "D2", "E2"
foreach ($col in $dateCols) {
    $cell = $ws.Range($col)
    $cell.NumberFormat = "@"
}

$ws.Range("A2").Value = "2024-07-03"
$ws.Range("D2").Value = "2024-07-08"
$ws.Range("E2").Value = "2024-07-15"

foreach ($col in $dateCols) {
    $ws.Range($col).Style = "Normal"
}

# Remaining columns are plain text/numbers and need no special formatting.
$ws.Range("B2").Value = "엑셀세라퓨틱스"
$ws.Range("C2").Value = "대신"
$ws.Range("F2").Value = 16180000
$ws.Range("G2").Value = 1618000
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = 6200
$ws.Range("J2").Value = 7700
$ws.Range("K2").Value = "-"
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = "-"
$ws.Range("N2").Value = "-"
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = "-"
$ws.Range("Q2").Value = "-"
$ws.Range("R2").Value = "517.7 : 1"
$ws.Range("S2").Value = "-"
$ws.Range("T2").Value = "-"
